# Auto commit at 2025-11-16  9:06:21.57
#
# Daily metrics refresh: the "Metrics" sheet holds the raw numbers (col B)
# that everything else (the "today" dashboard sheet, via formulas such as
# =Metrics!B2) derives from. Write the new day's figures into Metrics!B2:B13
# - the dependent formulas on "today" (B11:B22, E11:E22, F11:F22) and the
# TODAY()-1 date cell recalculate automatically. Finally restore the
# selections recorded for each sheet.

$wb = $excel.ActiveWorkbook

$metrics = $wb.Worksheets.Item("Metrics")
$today   = $wb.Worksheets.Item("today")

$metrics.Range("B2").Value  = 200374.95000000004
$metrics.Range("B3").Value  = 176015.29
$metrics.Range("B4").Value  = 62188.41
$metrics.Range("B5").Value  = 8228
$metrics.Range("B6").Value  = 4996620.7000000011
$metrics.Range("B7").Value  = 4218091.9700000007
$metrics.Range("B8").Value  = 1469148.24
$metrics.Range("B9").Value  = 194435
$metrics.Range("B10").Value = 33462001.690000009
$metrics.Range("B11").Value = 31493367.129999999
$metrics.Range("B12").Value = 11750870.280000001
$metrics.Range("B13").Value = 1292065

# Restore the recorded cursor/selection on the Metrics sheet ...
[void]$metrics.Activate()
[void]$metrics.Range("B2:B13").Select()

# ... then leave the "today" sheet selected/active (it was the active tab
# before the edit and stays that way), with its own recorded selection.
[void]$today.Activate()
[void]$today.Range("E6").Select()
